# Apply odds/statistics updates per the 2025-05-30 FlashScore refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").Value = 1.13
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 1.8
$ws.Range("M2").Value = 1.8
$ws.Range("N2").Value = 3.6
$ws.Range("O2").Value = 1.29
$ws.Range("P2").Value = 1.8
$ws.Range("Q2").Value = 2

# Row 5
$ws.Range("G5").Value = 6.4
$ws.Range("I5").Value = 1.39
$ws.Range("N5").Value = 1.47
$ws.Range("O5").Value = 2.32
$ws.Range("R5").Value = 1.65
$ws.Range("S5").Value = 1.98
$ws.Range("V5").Value = 20
$ws.Range("X5").Value = 60
$ws.Range("Y5").Value = 50
$ws.Range("Z5").Value = 16.5
$ws.Range("AB5").Value = 16.5
$ws.Range("AC5").Value = 60
$ws.Range("AD5").Value = 9
$ws.Range("AE5").Value = 7.9
$ws.Range("AG5").Value = 10
$ws.Range("AI5").Value = 21
$ws.Range("AJ5").Value = 400

# Row 9
$ws.Range("G9").Value = 2.32
$ws.Range("I9").Value = 3.2
$ws.Range("Q9").Value = 2.5
$ws.Range("R9").Value = 1.95
$ws.Range("S9").Value = 1.75
$ws.Range("V9").Value = 10
$ws.Range("W9").Value = 25
$ws.Range("X9").Value = 23
$ws.Range("Y9").Value = 40
$ws.Range("AE9").Value = 16.5
$ws.Range("AF9").Value = 12.5
$ws.Range("AH9").Value = 35
$ws.Range("AJ9").Value = 1250

# Row 10
$ws.Range("G10").Value = 2.3
$ws.Range("H10").Value = 2.77
$ws.Range("I10").Value = 3.6
$ws.Range("J10").Value = 1.09
$ws.Range("K10").Value = 6.5
$ws.Range("L10").Value = 1.35
$ws.Range("M10").Value = 3
$ws.Range("N10").Value = 2.05
$ws.Range("O10").Value = 1.72
$ws.Range("P10").Value = 1.42
$ws.Range("Q10").Value = 2.75
$ws.Range("R10").Value = 1.7
$ws.Range("S10").Value = 2.05
$ws.Range("T10").Value = 7.4
$ws.Range("W10").Value = 28
$ws.Range("X10").Value = 21
$ws.Range("Y10").Value = 29
$ws.Range("Z10").Value = 6.5
$ws.Range("AA10").Value = 5.6
$ws.Range("AC10").Value = 60
$ws.Range("AD10").Value = 9.75
$ws.Range("AE10").Value = 22
$ws.Range("AG10").Value = 65
$ws.Range("AI10").Value = 40
$ws.Range("AJ10").Value = 500

# Row 12
$ws.Range("N12").Value = 2.03
$ws.Range("O12").Value = 1.78
$ws.Range("T12").Value = 6.5
$ws.Range("AB12").Value = 17

# Row 14
$ws.Range("J14").Value = 1.08
$ws.Range("K14").Value = 8
$ws.Range("N14").Value = 2.4
$ws.Range("O14").Value = 1.53

# Row 16
$ws.Range("N16").Value = 1.83
$ws.Range("O16").Value = 1.98

# Row 17
$ws.Range("K17").Value = 9

# Row 18
$ws.Range("L18").Value = 1.33
$ws.Range("M18").Value = 3.25
$ws.Range("N18").Value = 2.08
$ws.Range("O18").Value = 1.73

# Row 20
$ws.Range("H20").Value = 3.7
$ws.Range("I20").Value = 4.6

# Row 21
$ws.Range("H21").Value = 3.55
$ws.Range("I21").Value = 3.95
$ws.Range("L21").Value = 1.26
$ws.Range("M21").Value = 3.2
$ws.Range("N21").Value = 1.78
$ws.Range("O21").Value = 1.83
$ws.Range("R21").Value = 1.7
$ws.Range("S21").Value = 1.9
$ws.Range("U21").Value = 8.5
$ws.Range("W21").Value = 14.5
$ws.Range("Z21").Value = 11
$ws.Range("AA21").Value = 7
$ws.Range("AD21").Value = 12
$ws.Range("AF21").Value = 13
$ws.Range("AG21").Value = 60
$ws.Range("AH21").Value = 35
$ws.Range("AI21").Value = 40

# Row 23
$ws.Range("G23").Value = 5.7
$ws.Range("V23").Value = 18.5
$ws.Range("Z23").Value = 11.5
$ws.Range("AA23").Value = 7.9
$ws.Range("AB23").Value = 18
$ws.Range("AE23").Value = 7

# Row 26
$ws.Range("G26").Value = 2.8
$ws.Range("I26").Value = 2.38
$ws.Range("L26").Value = 1.3
$ws.Range("M26").Value = 3.4
$ws.Range("N26").Value = 1.98
$ws.Range("O26").Value = 1.83
$ws.Range("R26").Value = 1.8
$ws.Range("S26").Value = 1.91
$ws.Range("V26").Value = 11
$ws.Range("X26").Value = 23
$ws.Range("Y26").Value = 34
$ws.Range("AD26").Value = 8
$ws.Range("AF26").Value = 9.5
$ws.Range("AJ26").Value = 251

# Row 27
$ws.Range("J27").Value = 1.1
$ws.Range("K27").Value = 7

# Row 32
$ws.Range("G32").Value = 2.4
$ws.Range("H32").Value = 3.8
$ws.Range("K32").Value = 17
$ws.Range("N32").Value = 1.6
$ws.Range("O32").Value = 2.3
$ws.Range("R32").Value = 1.53
$ws.Range("S32").Value = 2.38
$ws.Range("T32").Value = 11
$ws.Range("V32").Value = 9.5
$ws.Range("W32").Value = 23
$ws.Range("X32").Value = 17
$ws.Range("Z32").Value = 17
$ws.Range("AA32").Value = 7.5
$ws.Range("AB32").Value = 13
$ws.Range("AC32").Value = 41
$ws.Range("AJ32").Value = 126

# Row 33
$ws.Range("G33").Value = 2.9
$ws.Range("I33").Value = 2.35
$ws.Range("X33").Value = 21
$ws.Range("Y33").Value = 26
$ws.Range("AD33").Value = 11
$ws.Range("AE33").Value = 13
$ws.Range("AF33").Value = 9.5
$ws.Range("AG33").Value = 23
$ws.Range("AH33").Value = 17
$ws.Range("AI33").Value = 21

# Row 34
$ws.Range("G34").Value = 1.45
$ws.Range("I34").Value = 5.75
$ws.Range("P34").Value = 1.2
$ws.Range("Q34").Value = 4.33
$ws.Range("R34").Value = 1.5
$ws.Range("S34").Value = 2.5
$ws.Range("T34").Value = 12
$ws.Range("U34").Value = 10
$ws.Range("W34").Value = 12
$ws.Range("Z34").Value = 23
$ws.Range("AB34").Value = 15
$ws.Range("AE34").Value = 41
$ws.Range("AF34").Value = 19
$ws.Range("AG34").Value = 67
$ws.Range("AH34").Value = 41
$ws.Range("AI34").Value = 34

# Row 35
$ws.Range("N35").Value = 1.33
$ws.Range("O35").Value = 3.4

# Row 36
$ws.Range("G36").Value = 1.88
$ws.Range("H36").Value = 3.4
$ws.Range("L36").Value = 1.3
$ws.Range("M36").Value = 2.95
$ws.Range("N36").Value = 1.88
$ws.Range("O36").Value = 1.72
$ws.Range("R36").Value = 1.78
$ws.Range("S36").Value = 1.83
$ws.Range("T36").Value = 6.9
$ws.Range("V36").Value = 8.5
$ws.Range("W36").Value = 16
$ws.Range("X36").Value = 15.5
$ws.Range("Y36").Value = 28
$ws.Range("Z36").Value = 9.75
$ws.Range("AA36").Value = 6.6
$ws.Range("AC36").Value = 75
$ws.Range("AD36").Value = 10.75

# Row 37
$ws.Range("G37").Value = 2.3
$ws.Range("H37").Value = 3.05
$ws.Range("L37").Value = 1.33
$ws.Range("M37").Value = 2.8
$ws.Range("N37").Value = 1.98
$ws.Range("O37").Value = 1.65
$ws.Range("R37").Value = 1.72
$ws.Range("S37").Value = 1.88
$ws.Range("T37").Value = 7.3
$ws.Range("U37").Value = 11
$ws.Range("V37").Value = 9
$ws.Range("X37").Value = 19.5
$ws.Range("Y37").Value = 30
$ws.Range("Z37").Value = 8.5
$ws.Range("AA37").Value = 5.9
$ws.Range("AB37").Value = 13.5
$ws.Range("AC37").Value = 65
$ws.Range("AD37").Value = 9
$ws.Range("AH37").Value = 27
$ws.Range("AJ37").Value = 600

Write-Output "Applied 200 cell updates"
